$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75; this shifts the existing rows 75-150
# down to 76-151 (matching the dimension change A1:T150 -> A1:T151).
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new record.
$ws.Range("A75").Value = 10
$ws.Range("B75").Value = "Vega Modelo de Temuco"
$ws.Range("C75").Value = "La Araucanía"
$ws.Range("D75").Value = 44827
$ws.Range("D75").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E75").Value = 9
$ws.Range("F75").Value = "Fruta"
$ws.Range("G75").Value = 100107
$ws.Range("H75").Value = "Otros"
$ws.Range("I75").Value = 100107002
$ws.Range("J75").Value = "Chirimoya"
$ws.Range("K75").Value = "Cultivar IV Región"
$ws.Range("L75").Value = "Primera"
$ws.Range("M75").Value = 75
$ws.Range("N75").Value = 3600
$ws.Range("O75").Value = 3600
$ws.Range("P75").Value = 3600
$ws.Range("Q75").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R75").Value = "Provincia del Elquí"
$ws.Range("S75").Value = 3600
$ws.Range("T75").Value = 1
